$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# "Meetings" -> "Meeting semanal" (row 9)
$ws.Range("B9").Value = "Meeting semanal"

# Revert the last task row: clear its text and drop it back to the
# underlined blank-line style used for the other spacer rows (e.g. B15)
$ws.Range("B13").Value = ""
$ws.Range("B13").Font.Underline = $true
